$wb = $excel.ActiveWorkbook

# ----- Sheet 1: ACE_landing_page_data -----
# Add the new 2023 row at the top, shift existing years down, drop 2017.
$ws1 = $wb.Worksheets.Item("ACE_landing_page_data")

$ws1.Cells.Item(2,1).Value2 = 2023
$ws1.Cells.Item(2,2).Value2 = 459.75478594995815
$ws1.Cells.Item(2,3).Value2 = 9697453988.650034
$ws1.Cells.Item(2,4).Value2 = 21092665.666574594
$ws1.Cells.Item(2,5).Value2 = 0.96487719841975128
$ws1.Cells.Item(2,6).Value2 = 143.23688519223813
$ws1.Cells.Item(2,7).Value2 = 311.30388937283158
$ws1.Cells.Item(2,8).Value2 = -0.073401833037353037
$ws1.Cells.Item(2,9).Value2 = 0.02881493514775868
$ws1.Cells.Item(2,10).Value2 = 0.1103140194202783
$ws1.Cells.Item(2,11).Value2 = 0.087627141137734288
$ws1.Cells.Item(2,12).Value2 = 0.019373872608318132
$ws1.Cells.Item(2,13).Value2 = -0.078394576610962274
$ws1.Cells.Item(2,14).Value2 = 98.083295763940995
$ws1.Cells.Item(2,15).Value2 = 98.340013892985766
$ws1.Cells.Item(3,1).Value2 = 2022
$ws1.Cells.Item(3,2).Value2 = 496.17493574050184
$ws1.Cells.Item(3,3).Value2 = 9425848767.6962852
$ws1.Cells.Item(3,4).Value2 = 18997027.235220883
$ws1.Cells.Item(3,5).Value2 = 0.88713968411124977
$ws1.Cells.Item(3,6).Value2 = 140.5145737409685
$ws1.Cells.Item(3,7).Value2 = 337.78435051745646
$ws1.Cells.Item(3,8).Value2 = -0.3463303356472609
$ws1.Cells.Item(3,9).Value2 = 0.034652402249871983
$ws1.Cells.Item(3,10).Value2 = 0.58283680377669089
$ws1.Cells.Item(3,11).Value2 = 0.46942846730010523
$ws1.Cells.Item(3,12).Value2 = -0.022739401580627816
$ws1.Cells.Item(3,13).Value2 = -0.35153884325092255
$ws1.Cells.Item(3,14).Value2 = 95.336189642167525
$ws1.Cells.Item(3,15).Value2 = 88.569550751355436
$ws1.Cells.Item(4,1).Value2 = 2021
$ws1.Cells.Item(4,2).Value2 = 759.06067360768861
$ws1.Cells.Item(4,3).Value2 = 9110159844.213953
$ws1.Cells.Item(4,4).Value2 = 12001886.227243055
$ws1.Cells.Item(4,5).Value2 = 0.60373111305054561
$ws1.Cells.Item(4,6).Value2 = 143.78413901904744
$ws1.Cells.Item(4,7).Value2 = 520.90144028189241
$ws1.Cells.Item(4,8).Value2 = -0.25298124827727519
$ws1.Cells.Item(4,9).Value2 = -0.050058973750314206
$ws1.Cells.Item(4,10).Value2 = 0.27164281225738329
$ws1.Cells.Item(4,11).Value2 = 0.25218144054202773
$ws1.Cells.Item(4,12).Value2 = -0.082298725560155828
$ws1.Cells.Item(4,13).Value2 = -0.24633456259909303
$ws1.Cells.Item(4,14).Value2 = 92.143206196455068
$ws1.Cells.Item(4,15).Value2 = 55.956211366848507
$ws1.Cells.Item(5,1).Value2 = 2020
$ws1.Cells.Item(5,2).Value2 = 1016.1199727010782
$ws1.Cells.Item(5,3).Value2 = 9590237280.4976711
$ws1.Cells.Item(5,4).Value2 = 9438095.4396601785
$ws1.Cells.Item(5,5).Value2 = 0.4821434765789297
$ws1.Cells.Item(5,6).Value2 = 156.67858705635106
$ws1.Cells.Item(5,7).Value2 = 691.15739482266133
$ws1.Cells.Item(5,8).Value2 = 1.2101261251179785
$ws1.Cells.Item(5,9).Value2 = -0.043309607927320526
$ws1.Cells.Item(5,10).Value2 = -0.56713312367111612
$ws1.Cells.Item(5,11).Value2 = -0.50709822624514778
$ws1.Cells.Item(5,12).Value2 = 0.083028081427743317
$ws1.Cells.Item(5,13).Value2 = 1.2162327457674529
$ws1.Cells.Item(5,14).Value2 = 96.998870088000416
$ws1.Cells.Item(5,15).Value2 = 44.003088624797613
$ws1.Cells.Item(6,1).Value2 = 2019
$ws1.Cells.Item(6,2).Value2 = 459.75655468388112
$ws1.Cells.Item(6,3).Value2 = 10024389666.671915
$ws1.Cells.Item(6,4).Value2 = 21803690.593523942
$ws1.Cells.Item(6,5).Value2 = 0.97817354745151885
$ws1.Cells.Item(6,6).Value2 = 144.66715105835806
$ws1.Cells.Item(6,7).Value2 = 311.86137653756327
$ws1.Cells.Item(6,8).Value2 = -0.0026066784032476864
$ws1.Cells.Item(6,9).Value2 = 0.013900326498015536
$ws1.Cells.Item(6,10).Value2 = 0.016550145808913896
$ws1.Cells.Item(6,11).Value2 = 0.010367650232698367
$ws1.Cells.Item(6,12).Value2 = 0.005753828532776728
$ws1.Cells.Item(6,13).Value2 = -0.0016745752680884518
$ws1.Cells.Item(6,14).Value2 = 101.39003264980155
$ws1.Cells.Item(6,15).Value2 = 101.65501458089139
$ws1.Cells.Item(7,1).Value2 = 2018
$ws1.Cells.Item(7,2).Value2 = 460.95812427122047
$ws1.Cells.Item(7,3).Value2 = 9886957726.2055798
$ws1.Cells.Item(7,4).Value2 = 21448711.29419177
$ws1.Cells.Item(7,5).Value2 = 0.9681362494397312
$ws1.Cells.Item(7,6).Value2 = 143.8395231061688
$ws1.Cells.Item(7,7).Value2 = 312.38448787509338
$ws1.Cells.Item(7,8).Value2 = -0.035645247613669828
$ws1.Cells.Item(7,9).Value2 = 0.01566648765768841
$ws1.Cells.Item(7,10).Value2 = 0.053208360454890302
$ws1.Cells.Item(7,11).Value2 = 0.04924647508858393
$ws1.Cells.Item(7,12).Value2 = 0.0000041615133750116939
$ws1.Cells.Item(7,13).Value2 = -0.030183232751069267
$ws1.Cells.Item(7,14).Value2 = 100
$ws1.Cells.Item(7,15).Value2 = 100

# Update sheet1 selection to reflect the newly entered row (A2:O7, active cell A2)
$ws1.Range("A2:O7").Select() | Out-Null

# ----- Sheet 2: ANSP -----
# Rename "AirNav Ireland" -> "Airnav Ireland" and add new entry "IAA";
# re-sort the ANSP name list alphabetically (list grows from 38 to 39 names).
$ws2 = $wb.Worksheets.Item("ANSP")

$ws2.Cells.Item(2,1).Value2 = 'Airnav Ireland'
$ws2.Cells.Item(3,1).Value2 = 'Albcontrol'
$ws2.Cells.Item(4,1).Value2 = 'ANS CR'
$ws2.Cells.Item(5,1).Value2 = 'ARMATS'
$ws2.Cells.Item(6,1).Value2 = 'Austro Control'
$ws2.Cells.Item(7,1).Value2 = 'Avinor (Continental)'
$ws2.Cells.Item(8,1).Value2 = 'BHANSA'
$ws2.Cells.Item(9,1).Value2 = 'BULATSA'
$ws2.Cells.Item(10,1).Value2 = 'Croatia Control'
$ws2.Cells.Item(11,1).Value2 = 'DCAC Cyprus'
$ws2.Cells.Item(12,1).Value2 = 'DFS'
$ws2.Cells.Item(13,1).Value2 = 'DHMI'
$ws2.Cells.Item(14,1).Value2 = 'DSNA'
$ws2.Cells.Item(15,1).Value2 = 'EANS'
$ws2.Cells.Item(16,1).Value2 = 'ENAIRE'
$ws2.Cells.Item(17,1).Value2 = 'ENAV'
$ws2.Cells.Item(18,1).Value2 = 'Fintraffic ANS'
$ws2.Cells.Item(19,1).Value2 = 'HASP'
$ws2.Cells.Item(20,1).Value2 = 'HungaroControl'
$ws2.Cells.Item(21,1).Value2 = 'IAA'
$ws2.Cells.Item(22,1).Value2 = 'LFV'
$ws2.Cells.Item(23,1).Value2 = 'LGS'
$ws2.Cells.Item(24,1).Value2 = 'LPS'
$ws2.Cells.Item(25,1).Value2 = 'LVNL'
$ws2.Cells.Item(26,1).Value2 = 'MATS'
$ws2.Cells.Item(27,1).Value2 = 'M-NAV'
$ws2.Cells.Item(28,1).Value2 = 'MOLDATSA'
$ws2.Cells.Item(29,1).Value2 = 'MUAC'
$ws2.Cells.Item(30,1).Value2 = 'NATS (Continental)'
$ws2.Cells.Item(31,1).Value2 = 'NAV Portugal (Continental)'
$ws2.Cells.Item(32,1).Value2 = 'NAVIAIR'
$ws2.Cells.Item(33,1).Value2 = 'Oro Navigacija'
$ws2.Cells.Item(34,1).Value2 = 'PANSA'
$ws2.Cells.Item(35,1).Value2 = 'ROMATSA'
$ws2.Cells.Item(36,1).Value2 = 'Sakaeronavigatsia'
$ws2.Cells.Item(37,1).Value2 = 'skeyes'
$ws2.Cells.Item(38,1).Value2 = 'Skyguide'
$ws2.Cells.Item(39,1).Value2 = 'Slovenia Control'
$ws2.Cells.Item(40,1).Value2 = 'SMATSA'

# Activate the ANSP sheet/tab and select the newly appended last entry
$ws2.Activate() | Out-Null
$ws2.Range("A40").Select() | Out-Null
